# Insert a new data row above row 95 (pushes the existing rows 95:209 down to 96:210)
# and populate it with a new price observation for "Ajo" (Chino, Primera) at
# "Terminal Hortofruticola Agro Chillan".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(95).Insert()

$ws.Cells.Item(95, 1).Value  = 7
$ws.Cells.Item(95, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(95, 3).Value  = "Ñuble"
$ws.Cells.Item(95, 4).Value2 = 44638
$ws.Cells.Item(95, 5).Value  = 16
$ws.Cells.Item(95, 6).Value  = 100112003
$ws.Cells.Item(95, 7).Value  = "Ajo"
$ws.Cells.Item(95, 8).Value  = "Chino"
$ws.Cells.Item(95, 9).Value  = "Primera"
$ws.Cells.Item(95, 10).Value = 60
$ws.Cells.Item(95, 11).Value = 18000
$ws.Cells.Item(95, 12).Value = 19000
$ws.Cells.Item(95, 13).Value = 18500
$ws.Cells.Item(95, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(95, 15).Value = "China"
$ws.Cells.Item(95, 16).Value = 1850
$ws.Cells.Item(95, 17).Value = 10
$ws.Cells.Item(95, 18).Value = "Hortaliza"
